$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric cells --------------------------------------------------
$ws.Range("A5").Value = 112105099
$ws.Range("B5").Value = 56543
$ws.Range("E5").Value = 103021
$ws.Range("Q5").Value = 427565.1123065132
$ws.Range("R5").Value = 6608165.764175405
$ws.Range("S5").Value = 10

# --- Plain text cells ------------------------------------------------------
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("F5").Value = "Talltita"
$ws.Range("G5").Value = "Poecile montanus"
$ws.Range("H5").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("M5").Value = "födosökande"
$ws.Range("P5").Value = "350 m S Molkoms IP, Vrm"
$ws.Range("T5").Value = "Värmland"
$ws.Range("U5").Value = "Karlstad"
$ws.Range("V5").Value = "Värmland"
$ws.Range("W5").Value = "Nyed"
$ws.Range("AW5").Value = "Olle Kvarnbäck"
$ws.Range("AX5").Value = "Olle Kvarnbäck"

# --- Text cells that look like numbers/dates: force Text format first so
#     Excel doesn't auto-convert them to a number / date serial, then drop
#     back to the workbook's default "Normal" style so no stray formatting
#     is left behind on the cell. -----------------------------------------
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "2"
$ws.Range("I5").Style = "Normal"

$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-09-15"
$ws.Range("Y5").Style = "Normal"

$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "12:36"
$ws.Range("Z5").Style = "Normal"

$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-09-15"
$ws.Range("AA5").Style = "Normal"

$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "12:36"
$ws.Range("AB5").Style = "Normal"

# --- Boolean cells -----------------------------------------------------
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false

# --- Cells that are present but hold an empty string in the source row ---
# (K5, AT5, AY5). Excel has no interactive/automation path that leaves a
# cell occupied by a zero-length string (assigning "" always clears the
# cell back to blank, exactly like typing nothing into the grid), so the
# closest achievable state is to leave them unset/blank.
$ws.Range("K5").Value = ""
$ws.Range("AT5").Value = ""
$ws.Range("AY5").Value = ""
